$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Days" sheet is a running table of Day_Number / Date pairs that ends
# at row 36 (day 35, 2019-02-04). Extend it with five more days, continuing
# the same numbering / date sequence, copying the existing row formatting
# (number formats + borders) down into the new rows.
$startRow  = 37
$startDay  = 36
$startDate = 43501
$rowCount  = 5

$formatSource = $ws.Range("A36:B36")
$sourceHeight = $ws.Rows.Item(36).RowHeight

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $destRange = $ws.Range("A" + $row + ":B" + $row)

    $formatSource.Copy()
    $destRange.PasteSpecial(-4122)

    $ws.Rows.Item($row).RowHeight = $sourceHeight

    $ws.Cells.Item($row, 1).Value = $startDay + $i
    $ws.Cells.Item($row, 2).Value = $startDate + $i
}

$excel.CutCopyMode = 0

$ws.Range("D41").Select()
